$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 270,4
$data[0,0] = 2024
$data[0,1] = "Brasil"
$data[0,2] = "Banana (cacho)"
$data[0,3] = 2.27956351839145
$data[1,0] = 2024
$data[1,1] = "Brasil"
$data[1,2] = "Goiaba"
$data[1,3] = 2.486657992731841
$data[2,0] = 2024
$data[2,1] = "Brasil"
$data[2,2] = "Limão"
$data[2,3] = 1.484949447325685
$data[3,0] = 2024
$data[3,1] = "Brasil"
$data[3,2] = "Mamão"
$data[3,3] = 1.869350039154268
$data[4,0] = 2024
$data[4,1] = "Brasil"
$data[4,2] = "Manga"
$data[4,3] = 1.972365892100205
$data[5,0] = 2024
$data[5,1] = "Brasil"
$data[5,2] = "Maracujá"
$data[5,3] = 3.47163056537942
$data[6,0] = 2024
$data[6,1] = "Nordeste"
$data[6,2] = "Banana (cacho)"
$data[6,3] = 2.000750024871291
$data[7,0] = 2024
$data[7,1] = "Nordeste"
$data[7,2] = "Goiaba"
$data[7,3] = 2.468553253837003
$data[8,0] = 2024
$data[8,1] = "Nordeste"
$data[8,2] = "Limão"
$data[8,3] = 1.662281500364166
$data[9,0] = 2024
$data[9,1] = "Nordeste"
$data[9,2] = "Mamão"
$data[9,3] = 1.558023704098608
$data[10,0] = 2024
$data[10,1] = "Nordeste"
$data[10,2] = "Manga"
$data[10,3] = 1.91665325462285
$data[11,0] = 2024
$data[11,1] = "Nordeste"
$data[11,2] = "Maracujá"
$data[11,3] = 3.178003439643948
$data[12,0] = 2024
$data[12,1] = "Sergipe"
$data[12,2] = "Banana (cacho)"
$data[12,3] = 2.525275226369617
$data[13,0] = 2024
$data[13,1] = "Sergipe"
$data[13,2] = "Goiaba"
$data[13,3] = 2.348579040852576
$data[14,0] = 2024
$data[14,1] = "Sergipe"
$data[14,2] = "Limão"
$data[14,3] = 1.294869580238383
$data[15,0] = 2024
$data[15,1] = "Sergipe"
$data[15,2] = "Mamão"
$data[15,3] = 2.838939453897902
$data[16,0] = 2024
$data[16,1] = "Sergipe"
$data[16,2] = "Manga"
$data[16,3] = 2.718074718074718
$data[17,0] = 2024
$data[17,1] = "Sergipe"
$data[17,2] = "Maracujá"
$data[17,3] = 2.499192719800382
$data[18,0] = 2023
$data[18,1] = "Brasil"
$data[18,2] = "Banana (cacho)"
$data[18,3] = 2.038547906892512
$data[19,0] = 2023
$data[19,1] = "Brasil"
$data[19,2] = "Goiaba"
$data[19,3] = 2.26141126653656
$data[20,0] = 2023
$data[20,1] = "Brasil"
$data[20,2] = "Limão"
$data[20,3] = 1.426673489607525
$data[21,0] = 2023
$data[21,1] = "Brasil"
$data[21,2] = "Mamão"
$data[21,3] = 2.198240851848563
$data[22,0] = 2023
$data[22,1] = "Brasil"
$data[22,2] = "Manga"
$data[22,3] = 1.854036375140777
$data[23,0] = 2023
$data[23,1] = "Brasil"
$data[23,2] = "Maracujá"
$data[23,3] = 3.386842251092613
$data[24,0] = 2023
$data[24,1] = "Nordeste"
$data[24,2] = "Banana (cacho)"
$data[24,3] = 1.720735615825297
$data[25,0] = 2023
$data[25,1] = "Nordeste"
$data[25,2] = "Goiaba"
$data[25,3] = 2.190569690094637
$data[26,0] = 2023
$data[26,1] = "Nordeste"
$data[26,2] = "Limão"
$data[26,3] = 1.607351896622372
$data[27,0] = 2023
$data[27,1] = "Nordeste"
$data[27,2] = "Mamão"
$data[27,3] = 1.549804689160202
$data[28,0] = 2023
$data[28,1] = "Nordeste"
$data[28,2] = "Manga"
$data[28,3] = 1.906242457317593
$data[29,0] = 2023
$data[29,1] = "Nordeste"
$data[29,2] = "Maracujá"
$data[29,3] = 3.086038514281053
$data[30,0] = 2023
$data[30,1] = "Sergipe"
$data[30,2] = "Banana (cacho)"
$data[30,3] = 1.965603743439272
$data[31,0] = 2023
$data[31,1] = "Sergipe"
$data[31,2] = "Goiaba"
$data[31,3] = 2.012674448613381
$data[32,0] = 2023
$data[32,1] = "Sergipe"
$data[32,2] = "Limão"
$data[32,3] = 1.712337239878328
$data[33,0] = 2023
$data[33,1] = "Sergipe"
$data[33,2] = "Mamão"
$data[33,3] = 2.529772581246159
$data[34,0] = 2023
$data[34,1] = "Sergipe"
$data[34,2] = "Manga"
$data[34,3] = 2.387377280100289
$data[35,0] = 2023
$data[35,1] = "Sergipe"
$data[35,2] = "Maracujá"
$data[35,3] = 2.215980496838153
$data[36,0] = 2022
$data[36,1] = "Brasil"
$data[36,2] = "Banana (cacho)"
$data[36,3] = 1.651329601608699
$data[37,0] = 2022
$data[37,1] = "Brasil"
$data[37,2] = "Goiaba"
$data[37,3] = 1.827714036564227
$data[38,0] = 2022
$data[38,1] = "Brasil"
$data[38,2] = "Limão"
$data[38,3] = 1.18932594205514
$data[39,0] = 2022
$data[39,1] = "Brasil"
$data[39,2] = "Mamão"
$data[39,3] = 1.798894572581742
$data[40,0] = 2022
$data[40,1] = "Brasil"
$data[40,2] = "Manga"
$data[40,3] = 1.232855637515202
$data[41,0] = 2022
$data[41,1] = "Brasil"
$data[41,2] = "Maracujá"
$data[41,3] = 2.599292999017086
$data[42,0] = 2022
$data[42,1] = "Nordeste"
$data[42,2] = "Banana (cacho)"
$data[42,3] = 1.518537637398577
$data[43,0] = 2022
$data[43,1] = "Nordeste"
$data[43,2] = "Goiaba"
$data[43,3] = 1.792264330987263
$data[44,0] = 2022
$data[44,1] = "Nordeste"
$data[44,2] = "Limão"
$data[44,3] = 1.464982616597644
$data[45,0] = 2022
$data[45,1] = "Nordeste"
$data[45,2] = "Mamão"
$data[45,3] = 1.197247779535877
$data[46,0] = 2022
$data[46,1] = "Nordeste"
$data[46,2] = "Manga"
$data[46,3] = 1.198385580489361
$data[47,0] = 2022
$data[47,1] = "Nordeste"
$data[47,2] = "Maracujá"
$data[47,3] = 2.469636989329464
$data[48,0] = 2022
$data[48,1] = "Sergipe"
$data[48,2] = "Banana (cacho)"
$data[48,3] = 1.704359443890875
$data[49,0] = 2022
$data[49,1] = "Sergipe"
$data[49,2] = "Goiaba"
$data[49,3] = 1.624974789532233
$data[50,0] = 2022
$data[50,1] = "Sergipe"
$data[50,2] = "Limão"
$data[50,3] = 2.206249804745929
$data[51,0] = 2022
$data[51,1] = "Sergipe"
$data[51,2] = "Mamão"
$data[51,3] = 1.902284315911647
$data[52,0] = 2022
$data[52,1] = "Sergipe"
$data[52,2] = "Manga"
$data[52,3] = 1.751199190621819
$data[53,0] = 2022
$data[53,1] = "Sergipe"
$data[53,2] = "Maracujá"
$data[53,3] = 1.702341584617215
$data[54,0] = 2021
$data[54,1] = "Brasil"
$data[54,2] = "Banana (cacho)"
$data[54,3] = 1.534800946363992
$data[55,0] = 2021
$data[55,1] = "Brasil"
$data[55,2] = "Goiaba"
$data[55,3] = 1.838723811516457
$data[56,0] = 2021
$data[56,1] = "Brasil"
$data[56,2] = "Limão"
$data[56,3] = 1.041715720904194
$data[57,0] = 2021
$data[57,1] = "Brasil"
$data[57,2] = "Mamão"
$data[57,3] = 1.182858500461295
$data[58,0] = 2021
$data[58,1] = "Brasil"
$data[58,2] = "Manga"
$data[58,3] = 1.35560943527196
$data[59,0] = 2021
$data[59,1] = "Brasil"
$data[59,2] = "Maracujá"
$data[59,3] = 2.344573057862099
$data[60,0] = 2021
$data[60,1] = "Nordeste"
$data[60,2] = "Banana (cacho)"
$data[60,3] = 1.41608069204233
$data[61,0] = 2021
$data[61,1] = "Nordeste"
$data[61,2] = "Goiaba"
$data[61,3] = 1.830666782181373
$data[62,0] = 2021
$data[62,1] = "Nordeste"
$data[62,2] = "Limão"
$data[62,3] = 1.388139959383111
$data[63,0] = 2021
$data[63,1] = "Nordeste"
$data[63,2] = "Mamão"
$data[63,3] = 1.105737789287738
$data[64,0] = 2021
$data[64,1] = "Nordeste"
$data[64,2] = "Manga"
$data[64,3] = 1.342914703046207
$data[65,0] = 2021
$data[65,1] = "Nordeste"
$data[65,2] = "Maracujá"
$data[65,3] = 2.174173469880545
$data[66,0] = 2021
$data[66,1] = "Sergipe"
$data[66,2] = "Banana (cacho)"
$data[66,3] = 1.562069681002609
$data[67,0] = 2021
$data[67,1] = "Sergipe"
$data[67,2] = "Goiaba"
$data[67,3] = 1.782710486991847
$data[68,0] = 2021
$data[68,1] = "Sergipe"
$data[68,2] = "Limão"
$data[68,3] = 1.376235977518972
$data[69,0] = 2021
$data[69,1] = "Sergipe"
$data[69,2] = "Mamão"
$data[69,3] = 1.367377935878815
$data[70,0] = 2021
$data[70,1] = "Sergipe"
$data[70,2] = "Manga"
$data[70,3] = 1.580450965486383
$data[71,0] = 2021
$data[71,1] = "Sergipe"
$data[71,2] = "Maracujá"
$data[71,3] = 1.65018371044969
$data[72,0] = 2020
$data[72,1] = "Brasil"
$data[72,2] = "Banana (cacho)"
$data[72,3] = 1.838707640729277
$data[73,0] = 2020
$data[73,1] = "Brasil"
$data[73,2] = "Goiaba"
$data[73,3] = 2.518004632326168
$data[74,0] = 2020
$data[74,1] = "Brasil"
$data[74,2] = "Limão"
$data[74,3] = 1.572023704085504
$data[75,0] = 2020
$data[75,1] = "Brasil"
$data[75,2] = "Mamão"
$data[75,3] = 1.275283153609351
$data[76,0] = 2020
$data[76,1] = "Brasil"
$data[76,2] = "Manga"
$data[76,3] = 1.651769088626799
$data[77,0] = 2020
$data[77,1] = "Brasil"
$data[77,2] = "Maracujá"
$data[77,3] = 2.81081939704712
$data[78,0] = 2020
$data[78,1] = "Nordeste"
$data[78,2] = "Banana (cacho)"
$data[78,3] = 1.665640487085381
$data[79,0] = 2020
$data[79,1] = "Nordeste"
$data[79,2] = "Goiaba"
$data[79,3] = 2.713530041050617
$data[80,0] = 2020
$data[80,1] = "Nordeste"
$data[80,2] = "Limão"
$data[80,3] = 1.73902622482806
$data[81,0] = 2020
$data[81,1] = "Nordeste"
$data[81,2] = "Mamão"
$data[81,3] = 1.267126600069477
$data[82,0] = 2020
$data[82,1] = "Nordeste"
$data[82,2] = "Manga"
$data[82,3] = 1.700535675229772
$data[83,0] = 2020
$data[83,1] = "Nordeste"
$data[83,2] = "Maracujá"
$data[83,3] = 2.598684893174639
$data[84,0] = 2020
$data[84,1] = "Sergipe"
$data[84,2] = "Banana (cacho)"
$data[84,3] = 1.811409843432785
$data[85,0] = 2020
$data[85,1] = "Sergipe"
$data[85,2] = "Goiaba"
$data[85,3] = 2.085130324057778
$data[86,0] = 2020
$data[86,1] = "Sergipe"
$data[86,2] = "Limão"
$data[86,3] = 1.676739813932897
$data[87,0] = 2020
$data[87,1] = "Sergipe"
$data[87,2] = "Mamão"
$data[87,3] = 1.67864081750555
$data[88,0] = 2020
$data[88,1] = "Sergipe"
$data[88,2] = "Manga"
$data[88,3] = 1.838164144600731
$data[89,0] = 2020
$data[89,1] = "Sergipe"
$data[89,2] = "Maracujá"
$data[89,3] = 2.312388579058331
$data[90,0] = 2019
$data[90,1] = "Brasil"
$data[90,2] = "Banana (cacho)"
$data[90,3] = 1.843359020101782
$data[91,0] = 2019
$data[91,1] = "Brasil"
$data[91,2] = "Goiaba"
$data[91,3] = 2.647689279737636
$data[92,0] = 2019
$data[92,1] = "Brasil"
$data[92,2] = "Limão"
$data[92,3] = 1.735543556788343
$data[93,0] = 2019
$data[93,1] = "Brasil"
$data[93,2] = "Mamão"
$data[93,3] = 1.547056506200332
$data[94,0] = 2019
$data[94,1] = "Brasil"
$data[94,2] = "Manga"
$data[94,3] = 1.933137845566006
$data[95,0] = 2019
$data[95,1] = "Brasil"
$data[95,2] = "Maracujá"
$data[95,3] = 3.323704087118848
$data[96,0] = 2019
$data[96,1] = "Nordeste"
$data[96,2] = "Banana (cacho)"
$data[96,3] = 1.718165947423648
$data[97,0] = 2019
$data[97,1] = "Nordeste"
$data[97,2] = "Goiaba"
$data[97,3] = 2.834150394982295
$data[98,0] = 2019
$data[98,1] = "Nordeste"
$data[98,2] = "Limão"
$data[98,3] = 1.88599859796241
$data[99,0] = 2019
$data[99,1] = "Nordeste"
$data[99,2] = "Mamão"
$data[99,3] = 1.286442902395676
$data[100,0] = 2019
$data[100,1] = "Nordeste"
$data[100,2] = "Manga"
$data[100,3] = 2.034693824302827
$data[101,0] = 2019
$data[101,1] = "Nordeste"
$data[101,2] = "Maracujá"
$data[101,3] = 3.084605105727647
$data[102,0] = 2019
$data[102,1] = "Sergipe"
$data[102,2] = "Banana (cacho)"
$data[102,3] = 2.080644725470809
$data[103,0] = 2019
$data[103,1] = "Sergipe"
$data[103,2] = "Goiaba"
$data[103,3] = 2.18605229508627
$data[104,0] = 2019
$data[104,1] = "Sergipe"
$data[104,2] = "Limão"
$data[104,3] = 2.61516510371253
$data[105,0] = 2019
$data[105,1] = "Sergipe"
$data[105,2] = "Mamão"
$data[105,3] = 1.739469975293001
$data[106,0] = 2019
$data[106,1] = "Sergipe"
$data[106,2] = "Manga"
$data[106,3] = 1.95735349538683
$data[107,0] = 2019
$data[107,1] = "Sergipe"
$data[107,2] = "Maracujá"
$data[107,3] = 2.504009653960393
$data[108,0] = 2018
$data[108,1] = "Brasil"
$data[108,2] = "Banana (cacho)"
$data[108,3] = 1.84492501083786
$data[109,0] = 2018
$data[109,1] = "Brasil"
$data[109,2] = "Goiaba"
$data[109,3] = 2.472460819184825
$data[110,0] = 2018
$data[110,1] = "Brasil"
$data[110,2] = "Limão"
$data[110,3] = 1.855331420334648
$data[111,0] = 2018
$data[111,1] = "Brasil"
$data[111,2] = "Mamão"
$data[111,3] = 1.561430079267445
$data[112,0] = 2018
$data[112,1] = "Brasil"
$data[112,2] = "Manga"
$data[112,3] = 1.806980256437124
$data[113,0] = 2018
$data[113,1] = "Brasil"
$data[113,2] = "Maracujá"
$data[113,3] = 3.011362106539858
$data[114,0] = 2018
$data[114,1] = "Nordeste"
$data[114,2] = "Banana (cacho)"
$data[114,3] = 1.764365398657053
$data[115,0] = 2018
$data[115,1] = "Nordeste"
$data[115,2] = "Goiaba"
$data[115,3] = 2.517430907512401
$data[116,0] = 2018
$data[116,1] = "Nordeste"
$data[116,2] = "Limão"
$data[116,3] = 2.504474641293347
$data[117,0] = 2018
$data[117,1] = "Nordeste"
$data[117,2] = "Mamão"
$data[117,3] = 1.54225840597576
$data[118,0] = 2018
$data[118,1] = "Nordeste"
$data[118,2] = "Manga"
$data[118,3] = 1.941145778756125
$data[119,0] = 2018
$data[119,1] = "Nordeste"
$data[119,2] = "Maracujá"
$data[119,3] = 2.819105323689936
$data[120,0] = 2018
$data[120,1] = "Sergipe"
$data[120,2] = "Banana (cacho)"
$data[120,3] = 2.18372499935639
$data[121,0] = 2018
$data[121,1] = "Sergipe"
$data[121,2] = "Goiaba"
$data[121,3] = 2.009675398701652
$data[122,0] = 2018
$data[122,1] = "Sergipe"
$data[122,2] = "Limão"
$data[122,3] = 2.261736257182142
$data[123,0] = 2018
$data[123,1] = "Sergipe"
$data[123,2] = "Mamão"
$data[123,3] = 1.520300575673253
$data[124,0] = 2018
$data[124,1] = "Sergipe"
$data[124,2] = "Manga"
$data[124,3] = 1.486637315331456
$data[125,0] = 2018
$data[125,1] = "Sergipe"
$data[125,2] = "Maracujá"
$data[125,3] = 1.660917359115673
$data[126,0] = 2017
$data[126,1] = "Brasil"
$data[126,2] = "Banana (cacho)"
$data[126,3] = 2.299701632469955
$data[127,0] = 2017
$data[127,1] = "Brasil"
$data[127,2] = "Goiaba"
$data[127,3] = 2.415799616763466
$data[128,0] = 2017
$data[128,1] = "Brasil"
$data[128,2] = "Limão"
$data[128,3] = 1.820819371328306
$data[129,0] = 2017
$data[129,1] = "Brasil"
$data[129,2] = "Mamão"
$data[129,3] = 1.67705410307074
$data[130,0] = 2017
$data[130,1] = "Brasil"
$data[130,2] = "Manga"
$data[130,3] = 1.751564886777807
$data[131,0] = 2017
$data[131,1] = "Brasil"
$data[131,2] = "Maracujá"
$data[131,3] = 2.980087114676682
$data[132,0] = 2017
$data[132,1] = "Nordeste"
$data[132,2] = "Banana (cacho)"
$data[132,3] = 2.186823882824749
$data[133,0] = 2017
$data[133,1] = "Nordeste"
$data[133,2] = "Goiaba"
$data[133,3] = 2.617391300331335
$data[134,0] = 2017
$data[134,1] = "Nordeste"
$data[134,2] = "Limão"
$data[134,3] = 1.719082790628175
$data[135,0] = 2017
$data[135,1] = "Nordeste"
$data[135,2] = "Mamão"
$data[135,3] = 1.664488397217599
$data[136,0] = 2017
$data[136,1] = "Nordeste"
$data[136,2] = "Manga"
$data[136,3] = 1.773282725006887
$data[137,0] = 2017
$data[137,1] = "Nordeste"
$data[137,2] = "Maracujá"
$data[137,3] = 2.569943897823992
$data[138,0] = 2017
$data[138,1] = "Sergipe"
$data[138,2] = "Banana (cacho)"
$data[138,3] = 2.8479927637762
$data[139,0] = 2017
$data[139,1] = "Sergipe"
$data[139,2] = "Goiaba"
$data[139,3] = 1.623003480072331
$data[140,0] = 2017
$data[140,1] = "Sergipe"
$data[140,2] = "Limão"
$data[140,3] = 2.61094592243772
$data[141,0] = 2017
$data[141,1] = "Sergipe"
$data[141,2] = "Mamão"
$data[141,3] = 1.622278588330398
$data[142,0] = 2017
$data[142,1] = "Sergipe"
$data[142,2] = "Manga"
$data[142,3] = 1.346032393091529
$data[143,0] = 2017
$data[143,1] = "Sergipe"
$data[143,2] = "Maracujá"
$data[143,3] = 1.864089782694982
$data[144,0] = 2016
$data[144,1] = "Brasil"
$data[144,2] = "Banana (cacho)"
$data[144,3] = 2.324076351451684
$data[145,0] = 2016
$data[145,1] = "Brasil"
$data[145,2] = "Goiaba"
$data[145,3] = 2.323555716465603
$data[146,0] = 2016
$data[146,1] = "Brasil"
$data[146,2] = "Limão"
$data[146,3] = 1.955353924566385
$data[147,0] = 2016
$data[147,1] = "Brasil"
$data[147,2] = "Mamão"
$data[147,3] = 2.029742836391957
$data[148,0] = 2016
$data[148,1] = "Brasil"
$data[148,2] = "Manga"
$data[148,3] = 1.636587015189683
$data[149,0] = 2016
$data[149,1] = "Brasil"
$data[149,2] = "Maracujá"
$data[149,3] = 2.933077365708656
$data[150,0] = 2016
$data[150,1] = "Nordeste"
$data[150,2] = "Banana (cacho)"
$data[150,3] = 2.248181807045193
$data[151,0] = 2016
$data[151,1] = "Nordeste"
$data[151,2] = "Goiaba"
$data[151,3] = 2.341769825804162
$data[152,0] = 2016
$data[152,1] = "Nordeste"
$data[152,2] = "Limão"
$data[152,3] = 1.521424394933396
$data[153,0] = 2016
$data[153,1] = "Nordeste"
$data[153,2] = "Mamão"
$data[153,3] = 1.647331402707658
$data[154,0] = 2016
$data[154,1] = "Nordeste"
$data[154,2] = "Manga"
$data[154,3] = 1.479186017704671
$data[155,0] = 2016
$data[155,1] = "Nordeste"
$data[155,2] = "Maracujá"
$data[155,3] = 2.448388698595438
$data[156,0] = 2016
$data[156,1] = "Sergipe"
$data[156,2] = "Banana (cacho)"
$data[156,3] = 2.486976249062884
$data[157,0] = 2016
$data[157,1] = "Sergipe"
$data[157,2] = "Goiaba"
$data[157,3] = 1.649334275372442
$data[158,0] = 2016
$data[158,1] = "Sergipe"
$data[158,2] = "Limão"
$data[158,3] = 2.823768357235594
$data[159,0] = 2016
$data[159,1] = "Sergipe"
$data[159,2] = "Mamão"
$data[159,3] = 2.424143842111441
$data[160,0] = 2016
$data[160,1] = "Sergipe"
$data[160,2] = "Manga"
$data[160,3] = 1.882923334347782
$data[161,0] = 2016
$data[161,1] = "Sergipe"
$data[161,2] = "Maracujá"
$data[161,3] = 2.334900889022499
$data[162,0] = 2015
$data[162,1] = "Brasil"
$data[162,2] = "Banana (cacho)"
$data[162,3] = 1.788052875170639
$data[163,0] = 2015
$data[163,1] = "Brasil"
$data[163,2] = "Goiaba"
$data[163,3] = 2.373887118602053
$data[164,0] = 2015
$data[164,1] = "Brasil"
$data[164,2] = "Limão"
$data[164,3] = 1.516043251675552
$data[165,0] = 2015
$data[165,1] = "Brasil"
$data[165,2] = "Mamão"
$data[165,3] = 1.685236416827466
$data[166,0] = 2015
$data[166,1] = "Brasil"
$data[166,2] = "Manga"
$data[166,3] = 1.818706120222627
$data[167,0] = 2015
$data[167,1] = "Brasil"
$data[167,2] = "Maracujá"
$data[167,3] = 2.814902658981434
$data[168,0] = 2015
$data[168,1] = "Nordeste"
$data[168,2] = "Banana (cacho)"
$data[168,3] = 1.695790082718525
$data[169,0] = 2015
$data[169,1] = "Nordeste"
$data[169,2] = "Goiaba"
$data[169,3] = 2.401288605011612
$data[170,0] = 2015
$data[170,1] = "Nordeste"
$data[170,2] = "Limão"
$data[170,3] = 1.287721681292117
$data[171,0] = 2015
$data[171,1] = "Nordeste"
$data[171,2] = "Mamão"
$data[171,3] = 1.609936124780118
$data[172,0] = 2015
$data[172,1] = "Nordeste"
$data[172,2] = "Manga"
$data[172,3] = 1.709174451954977
$data[173,0] = 2015
$data[173,1] = "Nordeste"
$data[173,2] = "Maracujá"
$data[173,3] = 2.377309275114619
$data[174,0] = 2015
$data[174,1] = "Sergipe"
$data[174,2] = "Banana (cacho)"
$data[174,3] = 1.980681749664956
$data[175,0] = 2015
$data[175,1] = "Sergipe"
$data[175,2] = "Goiaba"
$data[175,3] = 1.941196807839499
$data[176,0] = 2015
$data[176,1] = "Sergipe"
$data[176,2] = "Limão"
$data[176,3] = 1.777047609939765
$data[177,0] = 2015
$data[177,1] = "Sergipe"
$data[177,2] = "Mamão"
$data[177,3] = 1.968657225670221
$data[178,0] = 2015
$data[178,1] = "Sergipe"
$data[178,2] = "Manga"
$data[178,3] = 1.737972303138866
$data[179,0] = 2015
$data[179,1] = "Sergipe"
$data[179,2] = "Maracujá"
$data[179,3] = 2.033207073489167
$data[180,0] = 2014
$data[180,1] = "Brasil"
$data[180,2] = "Banana (cacho)"
$data[180,3] = 1.798052046544692
$data[181,0] = 2014
$data[181,1] = "Brasil"
$data[181,2] = "Goiaba"
$data[181,3] = 2.238511707743075
$data[182,0] = 2014
$data[182,1] = "Brasil"
$data[182,2] = "Limão"
$data[182,3] = 1.633660822452277
$data[183,0] = 2014
$data[183,1] = "Brasil"
$data[183,2] = "Mamão"
$data[183,3] = 1.692190171838571
$data[184,0] = 2014
$data[184,1] = "Brasil"
$data[184,2] = "Manga"
$data[184,3] = 1.58982547765587
$data[185,0] = 2014
$data[185,1] = "Brasil"
$data[185,2] = "Maracujá"
$data[185,3] = 2.680755611778694
$data[186,0] = 2014
$data[186,1] = "Nordeste"
$data[186,2] = "Banana (cacho)"
$data[186,3] = 1.651921776531765
$data[187,0] = 2014
$data[187,1] = "Nordeste"
$data[187,2] = "Goiaba"
$data[187,3] = 2.232817050949128
$data[188,0] = 2014
$data[188,1] = "Nordeste"
$data[188,2] = "Limão"
$data[188,3] = 1.575060983406571
$data[189,0] = 2014
$data[189,1] = "Nordeste"
$data[189,2] = "Mamão"
$data[189,3] = 1.672830123478558
$data[190,0] = 2014
$data[190,1] = "Nordeste"
$data[190,2] = "Manga"
$data[190,3] = 1.494643764323979
$data[191,0] = 2014
$data[191,1] = "Nordeste"
$data[191,2] = "Maracujá"
$data[191,3] = 2.217020590037997
$data[192,0] = 2014
$data[192,1] = "Sergipe"
$data[192,2] = "Banana (cacho)"
$data[192,3] = 2.177266257570078
$data[193,0] = 2014
$data[193,1] = "Sergipe"
$data[193,2] = "Goiaba"
$data[193,3] = 1.847236466850893
$data[194,0] = 2014
$data[194,1] = "Sergipe"
$data[194,2] = "Limão"
$data[194,3] = 1.794179666832193
$data[195,0] = 2014
$data[195,1] = "Sergipe"
$data[195,2] = "Mamão"
$data[195,3] = 2.140867358891868
$data[196,0] = 2014
$data[196,1] = "Sergipe"
$data[196,2] = "Manga"
$data[196,3] = 2.050441834470599
$data[197,0] = 2014
$data[197,1] = "Sergipe"
$data[197,2] = "Maracujá"
$data[197,3] = 2.15853301363992
$data[198,0] = 2013
$data[198,1] = "Brasil"
$data[198,2] = "Banana (cacho)"
$data[198,3] = 1.738082342401291
$data[199,0] = 2013
$data[199,1] = "Brasil"
$data[199,2] = "Goiaba"
$data[199,3] = 2.630629717287714
$data[200,0] = 2013
$data[200,1] = "Brasil"
$data[200,2] = "Limão"
$data[200,3] = 1.37530145432961
$data[201,0] = 2013
$data[201,1] = "Brasil"
$data[201,2] = "Mamão"
$data[201,3] = 1.790601247526734
$data[202,0] = 2013
$data[202,1] = "Brasil"
$data[202,2] = "Manga"
$data[202,3] = 1.826702255592424
$data[203,0] = 2013
$data[203,1] = "Brasil"
$data[203,2] = "Maracujá"
$data[203,3] = 2.604544551711273
$data[204,0] = 2013
$data[204,1] = "Nordeste"
$data[204,2] = "Banana (cacho)"
$data[204,3] = 1.681184049287278
$data[205,0] = 2013
$data[205,1] = "Nordeste"
$data[205,2] = "Goiaba"
$data[205,3] = 3.231497032509516
$data[206,0] = 2013
$data[206,1] = "Nordeste"
$data[206,2] = "Limão"
$data[206,3] = 1.232468795591848
$data[207,0] = 2013
$data[207,1] = "Nordeste"
$data[207,2] = "Mamão"
$data[207,3] = 1.798919733187058
$data[208,0] = 2013
$data[208,1] = "Nordeste"
$data[208,2] = "Manga"
$data[208,3] = 1.936650039990547
$data[209,0] = 2013
$data[209,1] = "Nordeste"
$data[209,2] = "Maracujá"
$data[209,3] = 2.273433436052047
$data[210,0] = 2013
$data[210,1] = "Sergipe"
$data[210,2] = "Banana (cacho)"
$data[210,3] = 2.535964233158265
$data[211,0] = 2013
$data[211,1] = "Sergipe"
$data[211,2] = "Goiaba"
$data[211,3] = 3.010833550922547
$data[212,0] = 2013
$data[212,1] = "Sergipe"
$data[212,2] = "Limão"
$data[212,3] = 1.784380662841208
$data[213,0] = 2013
$data[213,1] = "Sergipe"
$data[213,2] = "Mamão"
$data[213,3] = 2.28493569457894
$data[214,0] = 2013
$data[214,1] = "Sergipe"
$data[214,2] = "Manga"
$data[214,3] = 1.939332920295549
$data[215,0] = 2013
$data[215,1] = "Sergipe"
$data[215,2] = "Maracujá"
$data[215,3] = 2.682577184982722
$data[216,0] = 2012
$data[216,1] = "Brasil"
$data[216,2] = "Banana (cacho)"
$data[216,3] = 1.580463042729988
$data[217,0] = 2012
$data[217,1] = "Brasil"
$data[217,2] = "Goiaba"
$data[217,3] = 2.384619263952295
$data[218,0] = 2012
$data[218,1] = "Brasil"
$data[218,2] = "Limão"
$data[218,3] = 1.135797427533519
$data[219,0] = 2012
$data[219,1] = "Brasil"
$data[219,2] = "Mamão"
$data[219,3] = 1.906119016235578
$data[220,0] = 2012
$data[220,1] = "Brasil"
$data[220,2] = "Manga"
$data[220,3] = 1.393213961739158
$data[221,0] = 2012
$data[221,1] = "Brasil"
$data[221,2] = "Maracujá"
$data[221,3] = 2.741857954922593
$data[222,0] = 2012
$data[222,1] = "Nordeste"
$data[222,2] = "Banana (cacho)"
$data[222,3] = 1.46765214649091
$data[223,0] = 2012
$data[223,1] = "Nordeste"
$data[223,2] = "Goiaba"
$data[223,3] = 2.517100209158792
$data[224,0] = 2012
$data[224,1] = "Nordeste"
$data[224,2] = "Limão"
$data[224,3] = 1.264707712515186
$data[225,0] = 2012
$data[225,1] = "Nordeste"
$data[225,2] = "Mamão"
$data[225,3] = 1.980758506425043
$data[226,0] = 2012
$data[226,1] = "Nordeste"
$data[226,2] = "Manga"
$data[226,3] = 1.351171587443878
$data[227,0] = 2012
$data[227,1] = "Nordeste"
$data[227,2] = "Maracujá"
$data[227,3] = 2.45018292332202
$data[228,0] = 2012
$data[228,1] = "Sergipe"
$data[228,2] = "Banana (cacho)"
$data[228,3] = 2.034755231774039
$data[229,0] = 2012
$data[229,1] = "Sergipe"
$data[229,2] = "Goiaba"
$data[229,3] = 2.336761957396841
$data[230,0] = 2012
$data[230,1] = "Sergipe"
$data[230,2] = "Limão"
$data[230,3] = 2.011801386065533
$data[231,0] = 2012
$data[231,1] = "Sergipe"
$data[231,2] = "Mamão"
$data[231,3] = 2.150498606389584
$data[232,0] = 2012
$data[232,1] = "Sergipe"
$data[232,2] = "Manga"
$data[232,3] = 1.706247894226258
$data[233,0] = 2012
$data[233,1] = "Sergipe"
$data[233,2] = "Maracujá"
$data[233,3] = 1.854230886255815
$data[234,0] = 2011
$data[234,1] = "Brasil"
$data[234,2] = "Banana (cacho)"
$data[234,3] = 1.568030292497216
$data[235,0] = 2011
$data[235,1] = "Brasil"
$data[235,2] = "Goiaba"
$data[235,3] = 2.119628439935699
$data[236,0] = 2011
$data[236,1] = "Brasil"
$data[236,2] = "Limão"
$data[236,3] = 1.194934334683075
$data[237,0] = 2011
$data[237,1] = "Brasil"
$data[237,2] = "Mamão"
$data[237,3] = 1.831371384427347
$data[238,0] = 2011
$data[238,1] = "Brasil"
$data[238,2] = "Manga"
$data[238,3] = 1.369333684846055
$data[239,0] = 2011
$data[239,1] = "Brasil"
$data[239,2] = "Maracujá"
$data[239,3] = 2.423437751672251
$data[240,0] = 2011
$data[240,1] = "Nordeste"
$data[240,2] = "Banana (cacho)"
$data[240,3] = 1.394642179910411
$data[241,0] = 2011
$data[241,1] = "Nordeste"
$data[241,2] = "Goiaba"
$data[241,3] = 1.933873316621625
$data[242,0] = 2011
$data[242,1] = "Nordeste"
$data[242,2] = "Limão"
$data[242,3] = 1.140623587985496
$data[243,0] = 2011
$data[243,1] = "Nordeste"
$data[243,2] = "Mamão"
$data[243,3] = 2.004775919625166
$data[244,0] = 2011
$data[244,1] = "Nordeste"
$data[244,2] = "Manga"
$data[244,3] = 1.283914695360262
$data[245,0] = 2011
$data[245,1] = "Nordeste"
$data[245,2] = "Maracujá"
$data[245,3] = 2.207310477085146
$data[246,0] = 2011
$data[246,1] = "Sergipe"
$data[246,2] = "Banana (cacho)"
$data[246,3] = 1.738848267532301
$data[247,0] = 2011
$data[247,1] = "Sergipe"
$data[247,2] = "Goiaba"
$data[247,3] = 1.563973821081284
$data[248,0] = 2011
$data[248,1] = "Sergipe"
$data[248,2] = "Limão"
$data[248,3] = 1.210515437247683
$data[249,0] = 2011
$data[249,1] = "Sergipe"
$data[249,2] = "Mamão"
$data[249,3] = 2.067123205851634
$data[250,0] = 2011
$data[250,1] = "Sergipe"
$data[250,2] = "Manga"
$data[250,3] = 1.19740532501884
$data[251,0] = 2011
$data[251,1] = "Sergipe"
$data[251,2] = "Maracujá"
$data[251,3] = 1.490778844276103
$data[252,0] = 2010
$data[252,1] = "Brasil"
$data[252,2] = "Banana (cacho)"
$data[252,3] = 1.564893132013541
$data[253,0] = 2010
$data[253,1] = "Brasil"
$data[253,2] = "Goiaba"
$data[253,3] = 2.054892415051441
$data[254,0] = 2010
$data[254,1] = "Brasil"
$data[254,2] = "Limão"
$data[254,3] = 1.467178054874201
$data[255,0] = 2010
$data[255,1] = "Brasil"
$data[255,2] = "Mamão"
$data[255,3] = 2.280960969331264
$data[256,0] = 2010
$data[256,1] = "Brasil"
$data[256,2] = "Manga"
$data[256,3] = 1.450536225154542
$data[257,0] = 2010
$data[257,1] = "Brasil"
$data[257,2] = "Maracujá"
$data[257,3] = 2.485047065032755
$data[258,0] = 2010
$data[258,1] = "Nordeste"
$data[258,2] = "Banana (cacho)"
$data[258,3] = 1.434344085869846
$data[259,0] = 2010
$data[259,1] = "Nordeste"
$data[259,2] = "Goiaba"
$data[259,3] = 1.856017393187005
$data[260,0] = 2010
$data[260,1] = "Nordeste"
$data[260,2] = "Limão"
$data[260,3] = 1.279262826961926
$data[261,0] = 2010
$data[261,1] = "Nordeste"
$data[261,2] = "Mamão"
$data[261,3] = 2.344363544983891
$data[262,0] = 2010
$data[262,1] = "Nordeste"
$data[262,2] = "Manga"
$data[262,3] = 1.428788294335373
$data[263,0] = 2010
$data[263,1] = "Nordeste"
$data[263,2] = "Maracujá"
$data[263,3] = 2.351861352778077
$data[264,0] = 2010
$data[264,1] = "Sergipe"
$data[264,2] = "Banana (cacho)"
$data[264,3] = 1.795467445947326
$data[265,0] = 2010
$data[265,1] = "Sergipe"
$data[265,2] = "Goiaba"
$data[265,3] = 1.431939680277456
$data[266,0] = 2010
$data[266,1] = "Sergipe"
$data[266,2] = "Limão"
$data[266,3] = 1.386718809525041
$data[267,0] = 2010
$data[267,1] = "Sergipe"
$data[267,2] = "Mamão"
$data[267,3] = 2.126133790446892
$data[268,0] = 2010
$data[268,1] = "Sergipe"
$data[268,2] = "Manga"
$data[268,3] = 1.360045370757244
$data[269,0] = 2010
$data[269,1] = "Sergipe"
$data[269,2] = "Maracujá"
$data[269,3] = 1.496014341836193

$ws.Range("A2:D271").Value = $data
